# Update the embedded "generated on" timestamp throughout the deck:
#   08/11/2022 14:57:50  ->  09/11/2022 11:15:46
#
# Touches the subtitle's "Date : ..." line on slide 1, plus the
# " Worker Wanda -..." stamp textbox repeated on every slide.

$oldStamp = "08/11/2022 14:57:50"
$newStamp = "09/11/2022 11:15:46"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)

    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $shp = $s.Shapes.Item($j)

        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $paraCount = $tr.Paragraphs().Count

            for ($k = 1; $k -le $paraCount; $k++) {
                $para = $tr.Paragraphs($k, 1)
                $paraText = $para.Text

                if ($paraText -like "*$oldStamp*") {
                    $newText = $paraText.Replace($oldStamp, $newStamp)
                    $runCount = $para.Runs().Count

                    if ($runCount -eq 1) {
                        # Rewrite the single run in place so its
                        # formatting (rPr) is preserved untouched.
                        $para.Runs(1, 1).Text = $newText
                    } else {
                        $para.Text = $newText
                    }
                }
            }
        }
    }
}
